# Insert a new data row at row 203 (pushing existing rows 203-302 down to 204-303)
# and populate it with a new weekly price record for Mango / Vega Modelo de Temuco.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 203; all rows below shift down by one.
$ws.Rows.Item(203).Insert()

# Fill in the new row's data.
$ws.Range("A203").Value = 10
$ws.Range("B203").Value = "Vega Modelo de Temuco"
$ws.Range("C203").Value = "La Araucanía"
$ws.Range("D203").Value = 44609
$ws.Range("E203").Value = 9
$ws.Range("F203").Value = "Fruta"
$ws.Range("G203").Value = 100108
$ws.Range("H203").Value = "Tropicales y subtropicales"
$ws.Range("I203").Value = 100108002
$ws.Range("J203").Value = "Mango"
$ws.Range("K203").Value = "Sin especificar"
$ws.Range("L203").Value = "Primera"
$ws.Range("M203").Value = 1000
$ws.Range("N203").Value = 7000
$ws.Range("O203").Value = 8000
$ws.Range("P203").Value = 7600
$ws.Range("Q203").Value = '$/bandeja 4 kilos'
$ws.Range("R203").Value = "Perú"
$ws.Range("S203").Value = 1900
$ws.Range("T203").Value = 4
